# Adds the "3D Printed PLA" samara wing-loading study to the sheet:
#  - Wing Loading (N/m^2) column (E) for the existing PLA-wax samara set (rows 69-98)
#  - a new block (rows 99-108) with the 3D printed samara mass/length/area data
#    and its own wing-loading calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant used to (re)apply an existing cell's number/font format to new cells
# without disturbing its value - PasteSpecial xlPasteFormats.
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Row 68 - new bold header "Wing Loading (N/m^2)" above the E column.
# ---------------------------------------------------------------------------
$ws.Range("B68").Copy()
$ws.Range("E68").PasteSpecial($xlPasteFormats)
$ws.Range("E68").Value = "Wing Loading (N/m^2)"

# ---------------------------------------------------------------------------
# Rows 69-98 - Wing Loading (N/m^2) formula filled down column E.
# Pick up the plain-data style already used throughout this block (B69) so
# every cell in E69:E98 - including the ones that did not exist before -
# ends up with the same formatting as its row neighbours.
# ---------------------------------------------------------------------------
$ws.Range("B69").Copy()
$ws.Range("E69:E98").PasteSpecial($xlPasteFormats)
$ws.Range("E69:E98").Formula = "=W5*10^-6*9.81/(B69*10^-6)"

# ---------------------------------------------------------------------------
# Row 99 - section header.
# ---------------------------------------------------------------------------
$ws.Range("B68").Copy()
$ws.Range("D99").PasteSpecial($xlPasteFormats)
$ws.Range("D99").Value = "3D Printed PLA"

# ---------------------------------------------------------------------------
# Row 100 - "Mass mg" sample index header row + "Total Mass" label.
# ---------------------------------------------------------------------------
$ws.Range("B68").Copy()
$ws.Range("A100").PasteSpecial($xlPasteFormats)
$ws.Range("A100").Value = "Mass mg"

$ws.Range("B69").Copy()
$ws.Range("B100:U100").PasteSpecial($xlPasteFormats)
for ($c = 2; $c -le 21; $c++) {
    $ws.Cells.Item(100, $c).Value = $c - 1
}

$ws.Range("B69").Copy()
$ws.Range("W100").PasteSpecial($xlPasteFormats)
$ws.Range("W100").Value = "Total Mass"

# ---------------------------------------------------------------------------
# Row 101 - sample "1 (24)" mass measurements + running total.
# ---------------------------------------------------------------------------
$ws.Range("B69").Copy()
$ws.Range("A101:U101").PasteSpecial($xlPasteFormats)
$ws.Range("A101").Value = "1 (24)"
$row101 = @(22.6,21.5,31.0,33.9,16.9,5.9,5.9,5.7,5.9,5.2,3.4,3.9,3.6,3.8,4.3,3.3,3.4,2.0,1.7,1.3)
for ($c = 2; $c -le 21; $c++) {
    $ws.Cells.Item(101, $c).Value = $row101[$c - 2]
}
$ws.Range("W101:W102").Formula = "=SUM(B101:U101)"

# ---------------------------------------------------------------------------
# Row 102 - sample "2 (12)" mass measurements + running total.
# ---------------------------------------------------------------------------
$ws.Range("B69").Copy()
$ws.Range("A102:U102").PasteSpecial($xlPasteFormats)
$ws.Range("A102").Value = "2 (12)"
$row102 = @(27.0,28.3,30.6,24.7,14.2,6.2,5.7,5.1,5.8,5.7,4.0,4.9,4.2,3.7,3.2,3.3,2.6,1.0,0.9,0.8)
for ($c = 2; $c -le 21; $c++) {
    $ws.Cells.Item(102, $c).Value = $row102[$c - 2]
}

# ---------------------------------------------------------------------------
# Row 103 - section header.
# ---------------------------------------------------------------------------
$ws.Range("B68").Copy()
$ws.Range("A103").PasteSpecial($xlPasteFormats)
$ws.Range("A103").Value = "Length mm"

# ---------------------------------------------------------------------------
# Rows 104-105 - segment lengths (2.5 mm each) for the two printed samples.
# ---------------------------------------------------------------------------
$ws.Range("B69").Copy()
$ws.Range("A104:U105").PasteSpecial($xlPasteFormats)
$ws.Range("A104").Value = 1
$ws.Range("A105").Value = 2
for ($c = 2; $c -le 21; $c++) {
    $ws.Cells.Item(104, $c).Value = 2.5
    $ws.Cells.Item(105, $c).Value = 2.5
}

# ---------------------------------------------------------------------------
# Row 106 - section headers for the totals/ratio block.
# ---------------------------------------------------------------------------
$ws.Range("B68").Copy()
$ws.Range("A106").PasteSpecial($xlPasteFormats)
$ws.Range("E106").PasteSpecial($xlPasteFormats)
$ws.Range("A106").Value = "Total Area mm^2"
$ws.Range("E106").Value = "Wing Loading N/m^2"

# ---------------------------------------------------------------------------
# Rows 107-108 - per-sample total area + resulting wing loading.
# ---------------------------------------------------------------------------
$ws.Range("B69").Copy()
$ws.Range("A107:B108").PasteSpecial($xlPasteFormats)
$ws.Range("A107").Value = 1
$ws.Range("B107").Value = 554.577780214358
$ws.Range("A108").Value = 2
$ws.Range("B108").Value = 577.386823006657

$ws.Range("E107:E108").Formula = "=W101*10^-6*9.81/(B107*10^-6)"

# ---------------------------------------------------------------------------
# Cosmetic sheet-level outline defaults (best-effort; matches the authored
# workbook's <outlinePr summaryBelow="0" summaryRight="0"/>).
# ---------------------------------------------------------------------------
$ws.Outline.SummaryBelow = $false
$ws.Outline.SummaryRight = $false

Write-Output "Added 3D printed PLA samara data (rows 68-108)."
